$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 8 (ano = 2025) with refreshed metrics
$ws.Range("C8").Value = 1248
$ws.Range("E8").Value = 1047
$ws.Range("G8").Value = 83.89423076923077
$ws.Range("H8").Value = 16.10576923076923
